$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder comma-separated "IA Control" values in column A (and one Requirement text update in F192)
$ws.Range('A2').Value = 'AU-4,AU-4 (1)'
$ws.Range('A3').Value = 'CM-6 b,SC-5,SC-5 (2)'
$ws.Range('A4').Value = 'CM-5 (1),AU-7 b,AC-6 (9),AU-7 a,AC-6 (8),AU-12 (3),AU-8 b'
$ws.Range('A5').Value = 'AC-17 (9),CM-6 b,CM-7 b,AC-17 (1)'
$ws.Range('A8').Value = 'IA-2 (11),IA-2 (12)'
$ws.Range('A10').Value = 'CM-7 (5) (b),CM-7 (2)'
$ws.Range('A15').Value = 'IA-8,AU-3 (1),IA-2'
$ws.Range('A16').Value = 'CM-6 b,AC-6 (10)'
$ws.Range('A17').Value = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A19').Value = 'IA-5 (1) (a),CM-6 b,IA-5 (1) (b)'
$ws.Range('A21').Value = 'MA-4 (7),SC-10,MA-4 e,AC-12'
$ws.Range('A22').Value = 'CM-5 (1),AU-6 (4),AU-3 (1),AU-7 a,AU-12 a,AU-7 (1),CM-6 b,AU-14 (1),AU-3,MA-4 (1) (a)'
$ws.Range('A25').Value = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A29').Value = 'SC-8 (2),SC-8,SC-8 (1)'
$ws.Range('A31').Value = 'AU-12 c,AU-3 (1),AC-2 (4),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A34').Value = 'AC-11 b,AC-11 a'
$ws.Range('A39').Value = 'AU-3,CM-6 b'
$ws.Range('A45').Value = 'AC-8 c 1, AC-8 c 2, AC-8 c 3,AC-8 a,AC-8 b'
$ws.Range('A53').Value = 'MA-4 (6),SC-13'
$ws.Range('A56').Value = 'AU-12 c,MA-4 (1) (a)'
$ws.Range('A65').Value = 'CM-6 b,IA-2 (2)'
$ws.Range('A67').Value = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A69').Value = 'CM-5 (1),AU-12 c,AU-7 b,AU-7 a,AU-12 a,CM-6 b,AU-12 (3),AU-8 b'
$ws.Range('A71').Value = 'AU-3,AU-4 (1)'
$ws.Range('A77').Value = 'AU-12 c,AU-3 (1),AC-2 (4),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A79').Value = 'AU-9,AU-9 (3)'
$ws.Range('A80').Value = 'IA-2 (3),IA-2 (2),IA-2 (1),IA-2 (4)'
$ws.Range('A81').Value = 'CM-6 b,CM-5 (3)'
$ws.Range('A86').Value = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A88').Value = 'AU-12 c,CM-5 (1),AC-2 (4),AC-6 (9)'
$ws.Range('A89').Value = 'IA-2 (2),IA-2 (3),IA-2 (5),IA-2,IA-2 (4)'
$ws.Range('A90').Value = 'IA-2 (11),IA-2 (12)'
$ws.Range('A91').Value = 'AU-9,AU-9 (3)'
$ws.Range('A96').Value = 'SC-8,SC-8 (1),AC-18 (1)'
$ws.Range('A97').Value = 'AU-8 (1) (a),AU-8 (1) (b),AU-8 b'
$ws.Range('A102').Value = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A111').Value = 'AU-5 a,AU-5 b'
$ws.Range('A119').Value = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A124').Value = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A128').Value = 'CM-6 b,CM-7 a,IA-5 (1) (c)'
$ws.Range('A136').Value = 'AC-11 b,AC-11 (1)'
$ws.Range('A139').Value = 'CM-3 (5),SI-6 b,SI-6 d'
$ws.Range('A148').Value = 'AU-12 c,AU-3 (1),AU-12 a,AU-14 (1),AU-3,MA-4 (1) (a)'
$ws.Range('A157').Value = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A181').Value = 'CM-6 b,SC-3'
$ws.Range('F192').Value = 'Red Hat Enterprise Linux 9 must protect the confidentiality and integrity of all information at rest.'

# Fill in the previously-empty Fix (M42) cell with new guidance text
$m42Text = @"
Configure Red Hat Enterprise Linux 9 to prevent unauthorized modification of all information at rest by using disk encryption.
Encrypting a partition in an already installed system is more difficult, because existing partitions will need to be resized and changed.
To encrypt an entire partition, dedicate a partition for encryption in the partition layout.
"@
$ws.Range('M42').Value = $m42Text
